$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Train-Bridge Controller: PlantUML status "Complete" -> "Complete?" (Good -> Neutral style)
$ws.Range("H3").Value = "Complete?"
$ws.Range("H3").Style = "Neutral"

# Row 3 - Comments: add new comment
$ws.Range("I3").Value = "Turn train classes more generic? GPT sometimes creates duplicate classes."

# Row 7 - RCL Leader Election: PlantUML status "Not complete" -> "Changes done" (Bad -> Neutral style)
$ws.Range("H7").Value = "Changes done"
$ws.Range("H7").Style = "Neutral"

# Row 7 - Comments: update comment text
$ws.Range("I7").Value = "Now only one node diagram. How to tell LLM how to initiate generic classes?"

# Update last active selection to D8
[void]$ws.Range("D8").Select()
